# Refactored charges tables creation.
#
# The second row of the Charges sheet (the old blank spacer row, r=2) is
# removed; every row below it shifts up by one (row 3 -> 2, row 4 -> 3, ...,
# row 41 -> 40), so the data that used to start at row 3 now starts at row 2,
# and the table ends one row earlier (row 40 instead of row 41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank row 2 and shift everything below it up by one row.
$ws.Rows("2").Delete()

# Excel leaves the old selection behind after a row delete, so move the
# active selection onto the new row 2 (what used to be row 3's selection
# before the edit), matching what the author was last looking at.
$ws.Range("A2:XFD2").Select()
